# Append a fresh scrape run (2025-12-07 01:29 JST) to the "ランサーズ" sheet.
# A brand-new listing lands at the top (row 2, pushing the previously-seen
# rows down by one) and another new listing is appended at the very end
# (old last row moves from 8 -> 10). Every row's "取得日時" timestamp is
# refreshed to the new scrape time.
#
# Because the engine's Range.Insert()/Range.Hyperlinks scoping doesn't keep
# the per-row hyperlink relationships in sync with shifted rows, this
# rewrites the whole data block (rows 2-10) explicitly and then rebuilds
# the F-column hyperlinks from scratch in the correct order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-12-07 01:29:42"

# Columns: A=取得日時, B=タイトル, C=カテゴリ, D=価格, E=締切, F=URL, G=優先度スコア, H=スキル概要
$rows = @(
    @($timestamp, "【急募】あなたAIクローン構築パートナー募集!", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448719", 310, "🔥AI,Ai"),
    @($timestamp, "CapcutAPIを用いた動画の自動制作ツールの作成", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448662", 248, "🔥API ◆ツール"),
    @($timestamp, "CapcutAPIを用いた動画の自動制作ツールの作成", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448659", 248, "🔥API ◆ツール"),
    @($timestamp, "【自動化】Webサービス更新ツール開発(200アカウント管理)", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448409", 230, "◆ツール,開発 ◇管理"),
    @($timestamp, "【品質重視】出張買取サービス向け予約管理システム開発(UI/要件定義済/Cursor実装途中あり)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448677", 153, "◆開発,システム開発 ◇管理"),
    @($timestamp, "【急募】新規システム開発に伴う要件定義依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448563", 110, "◆開発,システム開発"),
    @($timestamp, "【受注メールを元にECサイト自動仕入ツール】", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448428", 98, "◆ツール ◇サイト"),
    @($timestamp, "【急募】HPリニューアルで業務フロー自動化を実現するプロ募集!", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448817", 88, "◆自動化"),
    @($timestamp, "【緊急】既存コードの構造解析ができるPHPエンジニアを探しています", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5448440", 33, "○PHP")
)

# Drop every existing hyperlink up front - they get rebuilt below once the
# URL text for each row is final, so there is no stale-target leftover.
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $data[5])
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
